$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

$ws.Range("C1").Value = "version"
$ws.Range("C2").Value = 1

$ws.Range("C3").Select()
